$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shared string text tweak: "Total Horas Diarias" -> "Total Horas Diarias  " (trailing spaces) ---
$ws.Range("D1").Value = "Total Horas Diarias  "

# --- New data rows (9, 10, 11) : Dia / Hora Inicio / Hora Fim ---
$ws.Range("A9").Value = 41845
$ws.Range("B9").Value = 0.95833333333333337
$ws.Range("C9").Value = 0.097222222222222224

$ws.Range("A10").Value = 41846
$ws.Range("B10").Value = 0.6875
$ws.Range("C10").Value = 0.83333333333333337

$ws.Range("A11").Value = 41849
$ws.Range("B11").Value = 0.95833333333333337
$ws.Range("C11").Value = 0.041666666666666664

# Match the date / time number formats used by the rest of column A / B / C
# (copy the existing style from row 2 so we reuse style indices 1 / 2
# instead of minting brand-new custom number formats)
$ws.Range("A2").Copy()
$ws.Range("A9:A11").PasteSpecial(-4122)
$ws.Range("B2:C2").Copy()
$ws.Range("B9:C11").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Column D ("Total Horas Diarias") formulas ---
# D3 keeps its own (non-shared) formula, set first so it is not folded into
# the later shared-formula groups.
$ws.Range("D3").Formula = "=+C3-B3"

# D2 is a standalone formula cell.
$ws.Range("D2").Formula = "=+C2-B2"

# D4:D11 share one formula pattern -> becomes a shared formula group.
$ws.Range("D4:D11").Formula = "=+C4-B4"

# --- Apply the custom time number format (h:mm;@) to the whole D2:D11 column ---
$ws.Range("D2:D11").NumberFormat = "h:mm;@"

# --- Remove the stray F17 (=25*400) cell/row ---
$ws.Range("F17").ClearContents()

# --- Selection / view state ---
$ws.Range("C12").Select()

# --- Window position tweak ---
$excel.Width = 1040
$excel.Top = 380
